$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from row 4 to row 5 so the date style (G) carries over
$ws.Range("A4:H4").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)

$ws.Range("A5").Value = 10027.950000000001
$ws.Range("B5").Value = 9987
$ws.Range("C5").Value = 80.11
$ws.Range("D5").Value = 79.78
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = -0.41
$ws.Range("G5").Value = 42609.505254629628
$ws.Range("H5").Value = $true
